$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 871.22784
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 856.7564
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 2570.2692
$ws.Range("M17").Value = -5832
$ws.Range("N17").Value = -2906.2692

$ws.Range("H121").Value = 4719.3335
$ws.Range("J121").Value = 4296.75
$ws.Range("L121").Value = 12890.25
$ws.Range("N121").Value = -16384.25

$ws.Range("H125").Value = 2115.5881
$ws.Range("J125").Value = 3599.625
$ws.Range("L125").Value = 32396.625
$ws.Range("N125").Value = -37316.625

$ws.Range("H132").Value = 5236.203
$ws.Range("I132").Value = 2908.862
$ws.Range("K132").Value = 8726.585999999999
$ws.Range("M132").Value = -6196.585999999999

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = $null
$ws.Range("N134").Value = 0

$ws.Range("H135").Value = 1061.8049
$ws.Range("I135").Value = 724.5
$ws.Range("J135").Value = 2261.111
$ws.Range("K135").Value = 6520.5
$ws.Range("L135").Value = 20349.999
$ws.Range("M135").Value = -3985.5
$ws.Range("N135").Value = -25419.999

$ws.Range("H137").Value = 4483.256
$ws.Range("I137").Value = 6671.136
$ws.Range("K137").Value = 20013.408
$ws.Range("M137").Value = -17463.408

$ws.Range("H139").Value = 69949
$ws.Range("J139").Value = 69949
$ws.Range("L139").Value = 69949
$ws.Range("N139").Value = -80229

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3440.4644
$ws.Range("I74").Value = 1255.6666
$ws.Range("K74").Value = 1255.6666
$ws.Range("M74").Value = -381.6666

$ws.Range("H77").Value = 3440.4644
$ws.Range("I77").Value = 1255.6666
$ws.Range("K77").Value = 6278.333000000001
$ws.Range("M77").Value = -1910.333000000001

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = $null
$ws.Range("N134").Value = 0

$ws.Range("H141").Value = 135064.64
$ws.Range("J141").Value = 135064.64
$ws.Range("L141").Value = 135064.64
$ws.Range("N141").Value = -145424.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3822.25
$ws.Range("I20").Value = 3308.6428
$ws.Range("K20").Value = 3308.6428
$ws.Range("M20").Value = -3061.6428

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = $null
$ws.Range("N108").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1652.1464
$ws.Range("I31").Value = 1323.7142
$ws.Range("J31").Value = 3568
$ws.Range("K31").Value = 1323.7142
$ws.Range("L31").Value = 3568
$ws.Range("M31").Value = -1028.7142
$ws.Range("N31").Value = -4158

$ws.Range("H34").Value = 1652.1464
$ws.Range("I34").Value = 1323.7142
$ws.Range("J34").Value = 3568
$ws.Range("K34").Value = 1323.7142
$ws.Range("L34").Value = 3568
$ws.Range("M34").Value = -1121.7142
$ws.Range("N34").Value = -3972

$ws.Range("H140").Value = 108516.6
$ws.Range("J140").Value = 112486
$ws.Range("L140").Value = 112486
$ws.Range("N140").Value = -122846

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1280.5555
$ws.Range("I4").Value = 168.27272
$ws.Range("J4").Value = 3028.4285
$ws.Range("K4").Value = 504.81816
$ws.Range("L4").Value = 9085.2855
$ws.Range("M4").Value = -392.81816
$ws.Range("N4").Value = -9309.2855

$ws.Range("H5").Value = 1762.4
$ws.Range("I5").Value = 1212.5454
$ws.Range("J5").Value = 3274.5
$ws.Range("K5").Value = 3637.6362
$ws.Range("L5").Value = 9823.5
$ws.Range("M5").Value = -3525.6362
$ws.Range("N5").Value = -10047.5

$ws.Range("H104").Value = 7161.8335
$ws.Range("J104").Value = 8394.200000000001
$ws.Range("L104").Value = 25182.6
$ws.Range("N104").Value = -30424.6

$ws.Range("H129").Value = 1355.625
$ws.Range("I129").Value = 824.3333
$ws.Range("J129").Value = 2949.5
$ws.Range("K129").Value = 2472.9999
$ws.Range("L129").Value = 8848.5
$ws.Range("M129").Value = 2527.0001
$ws.Range("N129").Value = -18848.5

$ws.Range("H134").Value = 3974.0833
$ws.Range("I134").Value = 3461.125
$ws.Range("K134").Value = 10383.375
$ws.Range("M134").Value = -5313.375

$ws.Range("H135").Value = 1762.4
$ws.Range("I135").Value = 1212.5454
$ws.Range("J135").Value = 3274.5
$ws.Range("K135").Value = 10912.9086
$ws.Range("L135").Value = 29470.5
$ws.Range("M135").Value = -8377.908599999999
$ws.Range("N135").Value = -34540.5

$ws.Range("H139").Value = 2811.9167
$ws.Range("I139").Value = 1451.7142
$ws.Range("K139").Value = 4355.142599999999
$ws.Range("M139").Value = 784.8574000000008

$ws.Range("H140").Value = 2029.8235
$ws.Range("I140").Value = 1304.9333
$ws.Range("K140").Value = 3914.7999
$ws.Range("M140").Value = 1265.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = $null
$ws.Range("N43").Value = 0

$ws.Range("H45").Value = 37057
$ws.Range("J45").Value = 37057
$ws.Range("L45").Value = 37057
$ws.Range("N45").Value = -38175

$ws.Range("H80").Value = 68574820
$ws.Range("I80").Value = 80003090
$ws.Range("K80").Value = 80003090
$ws.Range("M80").Value = -80002092

$ws.Range("H83").Value = 68574820
$ws.Range("I83").Value = 80003090
$ws.Range("K83").Value = 400015450
$ws.Range("M83").Value = -400010458

$ws.Range("H126").Value = 5116.357
$ws.Range("I126").Value = 5737.222
$ws.Range("J126").Value = 3998.8
$ws.Range("K126").Value = 17211.666
$ws.Range("L126").Value = 11996.4
$ws.Range("M126").Value = -14741.666
$ws.Range("N126").Value = -16936.4

$ws.Range("H132").Value = 4962.1934
$ws.Range("I132").Value = 3019.6458
$ws.Range("K132").Value = 9058.937399999999
$ws.Range("M132").Value = -6528.937399999999

$ws.Range("H134").Value = 255163
$ws.Range("J134").Value = 255163
$ws.Range("L134").Value = 765489
$ws.Range("N134").Value = -770559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 47621864
$ws.Range("I100").Value = 90911820
$ws.Range("J100").Value = 2899.5
$ws.Range("K100").Value = 90911820
$ws.Range("L100").Value = 2899.5
$ws.Range("M100").Value = -90911279
$ws.Range("N100").Value = -3981.5

$ws.Range("H136").Value = 5004503
$ws.Range("I136").Value = 7828284
$ws.Range("J136").Value = 8582
$ws.Range("K136").Value = 23484852
$ws.Range("L136").Value = 25746
$ws.Range("M136").Value = -23482302
$ws.Range("N136").Value = -30846

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4224.731
$ws.Range("J4").Value = 7295.4614
$ws.Range("L4").Value = 7295.4614
$ws.Range("N4").Value = -7521.4614

$ws.Range("H132").Value = 4701.8984
$ws.Range("I132").Value = 3508.0962
$ws.Range("K132").Value = 10524.2886
$ws.Range("M132").Value = -7994.2886

$ws.Range("H136").Value = 1724.2559
$ws.Range("I136").Value = 1851.1578
$ws.Range("K136").Value = 5553.4734
$ws.Range("M136").Value = -3003.4734

$ws.Range("H141").Value = 198000
$ws.Range("J141").Value = 198000
$ws.Range("L141").Value = 198000
$ws.Range("N141").Value = -208360
